$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("D3").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E3").Value = "['Normal', 'ParamViolation']"

# Row 25
$ws.Range("D25").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E25").Value = "['Normal']"

# Row 38
$ws.Range("D38").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E38").Value = "['Normal', 'SoftwareFault']"

# Row 39
$ws.Range("D39").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'SoftwareFault']"

# Row 56
$ws.Range("D56").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "['Normal']"

# Row 58
$ws.Range("D58").Value = "[0, 0, 0, 1, 0, 0, 1]"
$ws.Range("E58").Value = "['ParamViolation', 'SoftwareFault']"

# Row 74
$ws.Range("D74").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E74").Value = "['Normal']"

# Row 88
$ws.Range("D88").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E88").Value = "['Normal']"

# Row 113
$ws.Range("D113").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E113").Value = "['Normal', 'SoftwareFault']"

$wb.Save()
